$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (G33=5512)
$ws.Range("H33").Value = 213.45454
$ws.Range("I33").Value = 62.5
$ws.Range("J33").Value = 616
$ws.Range("K33").Value = 62.5
$ws.Range("L33").Value = 616
$ws.Range("M33").Value = 166.5
$ws.Range("N33").Value = -1074
# Row 39 (G39=4603)
$ws.Range("H39").Value = 242.66667
$ws.Range("I39").Value = 101.2
$ws.Range("J39").Value = 950
$ws.Range("K39").Value = 303.6
$ws.Range("L39").Value = 2850
$ws.Range("M39").Value = -7.600000000000023
$ws.Range("N39").Value = -3442
# Row 135 (G135=44047)
$ws.Range("H135").Value = 839.4545000000001
$ws.Range("I135").Value = 598.45
$ws.Range("K135").Value = 5386.05
$ws.Range("M135").Value = -2851.05
# Row 138 (G138=44169)
$ws.Range("H138").Value = 4118.4375
$ws.Range("J138").Value = 4532.6904
$ws.Range("L138").Value = 13598.0712
$ws.Range("N138").Value = -23878.0712

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 46 (G46=3498)
$ws.Range("H46").Value = 24999.75
$ws.Range("I46").Value = 26666.666
$ws.Range("K46").Value = 26666.666
$ws.Range("M46").Value = -26347.666
# Row 122 (G122=36168)
$ws.Range("H122").Value = 373049.28
$ws.Range("I122").Value = 557104.8
$ws.Range("K122").Value = 1671314.4
$ws.Range("M122").Value = -1668864.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94 (G94=19939)
$ws.Range("H94").Value = 4366.6665
$ws.Range("I94").Value = 4550
$ws.Range("K94").Value = 4550
$ws.Range("M94").Value = -4099
# Row 99 (G99=19943)
$ws.Range("H99").Value = 3642.25
$ws.Range("I99").Value = 3427.9092
$ws.Range("K99").Value = 3427.9092
$ws.Range("M99").Value = -1929.9092

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G31=44023)
$ws.Range("H31").Value = 5049.6665
$ws.Range("J31").Value = 6969
$ws.Range("L31").Value = 6969
$ws.Range("N31").Value = -7559
# Row 34 (G34=44023)
$ws.Range("H34").Value = 5049.6665
$ws.Range("J34").Value = 6969
$ws.Range("L34").Value = 6969
$ws.Range("N34").Value = -7373
# Row 44 (G44=1850)
$ws.Range("H44").Value = 17000
$ws.Range("J44").Value = 17000
$ws.Range("L44").Value = 17000
$ws.Range("N44").Value = -17884
# Row 92 (G92=18041)
$ws.Range("H92").Value = 28533.666
$ws.Range("J92").Value = 28533.666
$ws.Range("L92").Value = 28533.666
$ws.Range("N92").Value = -33525.666
# Row 122 (G122=36196)
$ws.Range("H122").Value = 3424.7104
$ws.Range("I122").Value = 3670.92
$ws.Range("J122").Value = 2951.2307
$ws.Range("K122").Value = 11012.76
$ws.Range("L122").Value = 8853.6921
$ws.Range("M122").Value = -8562.76
$ws.Range("N122").Value = -13753.6921

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8 (G8=16734)
$ws.Range("H8").Value = 196
$ws.Range("I8").Value = 196
$ws.Range("K8").Value = 588
$ws.Range("M8").Value = -449
# Row 38 (G38=4860)
$ws.Range("H38").Value = 67.90000000000001
$ws.Range("I38").Value = 36.75
$ws.Range("K38").Value = 110.25
$ws.Range("M38").Value = 236.75
# Row 107 (G107=27838)
$ws.Range("H107").Value = 708.6177
$ws.Range("J107").Value = 697.8125
$ws.Range("L107").Value = 2093.4375
$ws.Range("N107").Value = -5933.4375
# Row 139 (G139=44102)
$ws.Range("H139").Value = 5723.9287
$ws.Range("I139").Value = 2313.5715
$ws.Range("K139").Value = 6940.7145
$ws.Range("M139").Value = -1800.7145

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 10 (G10=4306)
$ws.Range("H10").Value = 3624996.5
$ws.Range("I10").Value = 4529996
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 4529996
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -4529827
$ws.Range("N10").Value = -5338
# Row 20 (G20=4095)
$ws.Range("H20").Value = 76060.71000000001
$ws.Range("J20").Value = 76060.71000000001
$ws.Range("L20").Value = 76060.71000000001
$ws.Range("N20").Value = -76550.71000000001
# Row 63 (G63=11048)
$ws.Range("H63").Value = 56722.25
$ws.Range("I63").Value = 59944.5
$ws.Range("J63").Value = 53500
$ws.Range("K63").Value = 59944.5
$ws.Range("L63").Value = 53500
$ws.Range("M63").Value = -59258.5
$ws.Range("N63").Value = -54872
# Row 66 (G66=11048)
$ws.Range("H66").Value = 56722.25
$ws.Range("I66").Value = 59944.5
$ws.Range("J66").Value = 53500
$ws.Range("K66").Value = 179833.5
$ws.Range("L66").Value = 160500
$ws.Range("M66").Value = -176401.5
$ws.Range("N66").Value = -167364
# Row 92 (G92=18094)
$ws.Range("H92").Value = 9250.200000000001
$ws.Range("J92").Value = 5312.75
$ws.Range("L92").Value = 5312.75
$ws.Range("N92").Value = -9056.75
# Row 102 (G102=36169)
$ws.Range("H102").Value = 1571.6774
$ws.Range("I102").Value = 351.63635
$ws.Range("J102").Value = 4554
$ws.Range("K102").Value = 351.63635
$ws.Range("L102").Value = 4554
$ws.Range("M102").Value = 1270.36365
$ws.Range("N102").Value = -7798

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G7=36249)
$ws.Range("H7").Value = 2297.4
$ws.Range("I7").Value = 2747
$ws.Range("J7").Value = 1997.6666
$ws.Range("K7").Value = 2747
$ws.Range("L7").Value = 1997.6666
$ws.Range("M7").Value = -2635
$ws.Range("N7").Value = -2221.6666
# Row 43 (G43=4314)
$ws.Range("H43").Value = 914546.75
$ws.Range("I43").Value = 500000
$ws.Range("J43").Value = 2020004.6
$ws.Range("K43").Value = 500000
$ws.Range("L43").Value = 2020004.6
$ws.Range("M43").Value = -499807
$ws.Range("N43").Value = -2020390.6
# Row 61 (G61=27740)
$ws.Range("H61").Value = 9595.6
$ws.Range("I61").Value = 9595.6
$ws.Range("K61").Value = 9595.6
$ws.Range("M61").Value = -9393.6
# Row 68 (G68=12563)
$ws.Range("H68").Value = 1672
$ws.Range("I68").Value = 1656.1
$ws.Range("J68").Value = 1725
$ws.Range("K68").Value = 1656.1
$ws.Range("L68").Value = 1725
$ws.Range("M68").Value = -907.0999999999999
$ws.Range("N68").Value = -3223
# Row 71 (G71=12563)
$ws.Range("H71").Value = 1672
$ws.Range("I71").Value = 1656.1
$ws.Range("J71").Value = 1725
$ws.Range("K71").Value = 8280.5
$ws.Range("L71").Value = 8625
$ws.Range("M71").Value = -4536.5
$ws.Range("N71").Value = -16113
# Row 82 (G82=12565)
$ws.Range("H82").Value = 4028.6667
$ws.Range("J82").Value = 1200
$ws.Range("L82").Value = 1200
$ws.Range("N82").Value = -1922
# Row 85 (G85=12565)
$ws.Range("H85").Value = 4028.6667
$ws.Range("J85").Value = 1200
$ws.Range("L85").Value = 1200
$ws.Range("N85").Value = -3696
# Row 93 (G93=19993)
$ws.Range("H93").Value = 674.375
$ws.Range("I93").Value = 565.8333
$ws.Range("K93").Value = 565.8333
$ws.Range("M93").Value = 682.1667
# Row 100 (G100=19995)
$ws.Range("H100").Value = 10333.333
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
# Row 113 (G113=27740)
$ws.Range("H113").Value = 9595.6
$ws.Range("I113").Value = 9595.6
$ws.Range("K113").Value = 9595.6
$ws.Range("M113").Value = -7425.6
# Row 122 (G122=36247)
$ws.Range("H122").Value = 4656.7896
$ws.Range("I122").Value = 2313.8462
$ws.Range("J122").Value = 9733.166999999999
$ws.Range("K122").Value = 6941.5386
$ws.Range("L122").Value = 29199.501
$ws.Range("M122").Value = -4491.5386
$ws.Range("N122").Value = -34099.501
# Row 126 (G126=36249)
$ws.Range("H126").Value = 2297.4
$ws.Range("I126").Value = 2747
$ws.Range("J126").Value = 1997.6666
$ws.Range("K126").Value = 8241
$ws.Range("L126").Value = 5992.9998
$ws.Range("M126").Value = -5771
$ws.Range("N126").Value = -10932.9998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 26 (G26=3800)
$ws.Range("H26").Value = 268626.5
$ws.Range("I26").Value = 24878
$ws.Range("K26").Value = 24878
$ws.Range("M26").Value = -24585
# Row 40 (G40=3601)
$ws.Range("H40").Value = 24588
$ws.Range("J40").Value = 24588
$ws.Range("L40").Value = 24588
$ws.Range("N40").Value = -24886
# Row 64 (G64=11036)
$ws.Range("H64").Value = 66495
$ws.Range("I64").Value = 64990
$ws.Range("J64").Value = 68000
$ws.Range("K64").Value = 64990
$ws.Range("L64").Value = 68000
$ws.Range("M64").Value = -64742
$ws.Range("N64").Value = -68496
# Row 67 (G67=11036)
$ws.Range("H67").Value = 66495
$ws.Range("I67").Value = 64990
$ws.Range("J67").Value = 68000
$ws.Range("K67").Value = 64990
$ws.Range("L67").Value = 68000
$ws.Range("M67").Value = -64132
$ws.Range("N67").Value = -69716
# Row 94 (G94=18075)
$ws.Range("H94").Value = 54745.832
$ws.Range("I94").Value = 69991.664
$ws.Range("K94").Value = 69991.664
$ws.Range("M94").Value = -69090.664
# Row 96 (G96=19977)
$ws.Range("H96").Value = 1212.5714
$ws.Range("I96").Value = 1249.6666
$ws.Range("K96").Value = 1249.6666
$ws.Range("M96").Value = 123.3334
# Row 100 (G100=19981)
$ws.Range("H100").Value = 3082.125
$ws.Range("I100").Value = 3159.5
$ws.Range("K100").Value = 6319
$ws.Range("M100").Value = -5778
# Row 113 (G113=27752)
$ws.Range("H113").Value = 1099.7059
$ws.Range("I113").Value = 1159.25
$ws.Range("J113").Value = 1046.7778
$ws.Range("K113").Value = 3477.75
$ws.Range("L113").Value = 3140.3334
$ws.Range("M113").Value = -1307.75
$ws.Range("N113").Value = -7480.3334
# Row 122 (G122=36208)
$ws.Range("H122").Value = 1814.5238
$ws.Range("I122").Value = 856.875
$ws.Range("K122").Value = 2570.625
$ws.Range("M122").Value = -120.625
# Row 126 (G126=36210)
$ws.Range("H126").Value = 2461.5557
$ws.Range("I126").Value = 894.25
$ws.Range("K126").Value = 2682.75
$ws.Range("M126").Value = -212.75
# Row 135 (G135=42043)
$ws.Range("H135").Value = 78997
$ws.Range("J135").Value = 78997
$ws.Range("L135").Value = 78997
$ws.Range("N135").Value = -89137
